$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTimestamp = "2023-01-28 12:56:19"

# Column O ("timestamp") holds this value for every data row (row 1 is the header).
# Determine the last used row dynamically and stamp the new timestamp on each data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 15).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 15).Value = $newTimestamp
}
